# "Domains and KPI's.xlsx" - Added Aerospace KPI sheet (manufacturing domain KPIs)
$wb = $excel.ActiveWorkbook

# --- The "Supply Chain Logistics" sheet was the active tab before this edit;
#     after adding the new sheet it is no longer active, and its lingering
#     selection/scroll position is reset to A1:C1 (matching the other,
#     previously-unselected sheets in the workbook). ---
$supplyChain = $wb.Worksheets.Item("Supply Chain Logistics")
$supplyChain.Range("A1:C1").Select()

# --- Add the new "Aerospace" worksheet as the last tab in the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Aerospace"

# Header row (Category / KPI / Short Description)
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "KPI"
$ws.Range("C1").Value = "Short Description"
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").VerticalAlignment = -4108
$ws.Range("B1:C1").WrapText = $true

# Row 2 - Safety
$ws.Range("B2").Value = "Satefy"
$ws.Range("C2").Value = " Building a safety culture in a plant should be intentional – from prevention to awareness to education to reporting. "
$ws.Range("C2").WrapText = $true

# Row 3 - First Yield Pass
$ws.Range("B3").Value = "First Yield Pass"
$ws.Range("C3").Value = "First Pass Yield is the percentage of products that pass all quality inspections and tests without requiring rework or repairs.igh FPY rates indicate that products are built right the first time, reducing costs associated with rework and increasing customer satisfaction."
$ws.Range("C3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 45

# Row 4 - On Time Delivery
$ws.Range("B4").Value = "On Time Delivery"
$ws.Range("C4").Value = "On-Time Delivery measures the percentage of products delivered to customers by the promised date."
$ws.Range("C4").WrapText = $true

# Category label for the Delivery group (rows 2-4)
$ws.Range("A2").Value = "Delivery"

# Row 5 - Cost Per Unit
$ws.Range("B5").Value = "Cost Per Unit"
$ws.Range("C5").Value = "Cost Per Unit tracks the total cost incurred to produce one unit of a product, including materials, labor, and overhead."
$ws.Range("C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 30

# Row 6 - Manufacturing Lead Time
$ws.Range("B6").Value = "Manufacturing Lead Time"
$ws.Range("C6").Value = "Manufacturing Lead Time measures the total time required to complete the production of a product, from the start of the manufacturing process to final delivery."
$ws.Range("C6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 30

# Row 7 - Supplier Quality Rating
$ws.Range("B7").Value = "Supplier Quality Rating"
$ws.Range("C7").Value = "Supplier Quality Rating assesses the performance and reliability of suppliers based on the quality of materials and components they provide."
$ws.Range("C7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 30

# Row 8 - ROI
$ws.Range("B8").Value = "ROI"
$ws.Range("C8").Value = "Return on Investment measures the profitability of investments in equipment, technology, or projects by comparing the return to the investment cost."
$ws.Range("C8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 30

# Row 9 - Compliance Rate
$ws.Range("B9").Value = "Compliance Rate"
$ws.Range("C9").Value = " Compliance Rate tracks the percentage of products and processes that adhere to industry regulations, standards, and contractual requirements."
$ws.Range("C9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30

# Row 10 - Inventory Turn Over Rate
$ws.Range("B10").Value = "Inventory Turn Over Rate"
$ws.Range("C10").Value = "Inventory Turnover Ratio measures how often inventory is sold and replaced over a given period."
$ws.Range("C10").WrapText = $true

# Row 11 - Scrap Rate
$ws.Range("B11").Value = "Scrap Rate"
$ws.Range("C11").Value = "Scrap Rate measures the percentage of materials or products that are discarded due to defects or errors during the manufacturing process."
$ws.Range("C11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 30

# Column widths (best-fit sized for KPI name / description columns)
$ws.Columns.Item(2).ColumnWidth = 22.83
$ws.Columns.Item(3).ColumnWidth = 105.3

# Page setup
$ws.PageSetup.Orientation = 1

# Final selection/view on the new, now-active sheet
$ws.Range("D8").Select()

Write-Output "Aerospace sheet added"
